$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the last-updated date (C1) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- "RAF-capacity" sheet: capacity credit multipliers updated for the
#     last two technologies (hydrogen combustion turbine / combined cycle) ---
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# Give column A a touch more breathing room to match the updated layout
$wsCapacity.Columns.Item(1).ColumnWidth = 28.166666666666668

# --- Make "RAF-capacity" the active/visible sheet, scrolled & zoomed the
#     way the author left it, with B25 selected ---
$wsCapacity.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 80
$win.ScrollRow = 14
$win.ScrollColumn = 1
$wsCapacity.Range("B25").Select() | Out-Null
